# Update countries & provincias Spain
# - Refresh COVID-19 stats for several countries (rows identified by current
#   country name in column A at the time each value is written).
# - "Republica del Chad" overtakes "Sierra Leona" and "Tanzania" in the
#   ranking (column A is sorted descending by total cases), so it moves up
#   three rows while the other two slide down one row each.
# - Update the "last updated" timestamp string in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" banner -----------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 21:05"

# --- Refresh stats for Estados Unidos (row 4) ------------------------------
$ws.Cells.Item(4, 2).Value = 1539030
$ws.Cells.Item(4, 3).Value = 11366
$ws.Cells.Item(4, 5).Value = 1096197
$ws.Cells.Item(4, 7).Value = 507
$ws.Cells.Item(4, 8).Value = 91485

# --- Refresh stats for Francia (row 10) ------------------------------------
$ws.Cells.Item(10, 4).Value = 61728
$ws.Cells.Item(10, 5).Value = 89602
$ws.Cells.Item(10, 7).Value = 131
$ws.Cells.Item(10, 8).Value = 28239

# --- Refresh stats for Costa Rica (row 113) --------------------------------
$ws.Cells.Item(113, 2).Value = 866
$ws.Cells.Item(113, 3).Value = 3
$ws.Cells.Item(113, 4).Value = 575
$ws.Cells.Item(113, 5).Value = 281

# --- Republica del Chad moves above Sierra Leona and Tanzania --------------
# Row 127 was "Sierra Leona", row 128 "Tanzania", row 129 "Republica del Chad".
# After the refresh, the updated totals re-sort the table so that
# "Republica del Chad" lands on row 127 (with its own new figures), and
# "Sierra Leona" / "Tanzania" shift down to rows 128 / 129 (keeping their
# previous figures).
$sierraLeona = @($ws.Cells.Item(127, 2).Value2, $ws.Cells.Item(127, 3).Value2, $ws.Cells.Item(127, 4).Value2, $ws.Cells.Item(127, 5).Value2, $ws.Cells.Item(127, 6).Value2, $ws.Cells.Item(127, 7).Value2, $ws.Cells.Item(127, 8).Value2)
$tanzania = @($ws.Cells.Item(128, 2).Value2, $ws.Cells.Item(128, 3).Value2, $ws.Cells.Item(128, 4).Value2, $ws.Cells.Item(128, 5).Value2, $ws.Cells.Item(128, 6).Value2, $ws.Cells.Item(128, 7).Value2, $ws.Cells.Item(128, 8).Value2)

# Row 127 -> Republica del Chad, with new stats
$ws.Cells.Item(127, 1).Value = "Republica del Chad"
$ws.Cells.Item(127, 2).Value = 519
$ws.Cells.Item(127, 3).Value = 16
$ws.Cells.Item(127, 4).Value = 117
$ws.Cells.Item(127, 5).Value = 349
$ws.Cells.Item(127, 6).Value = 0
$ws.Cells.Item(127, 7).Value = 0
$ws.Cells.Item(127, 8).Value = 53

# Row 128 -> Sierra Leona, keeping its previous figures
$ws.Cells.Item(128, 1).Value = "Sierra Leona"
$ws.Cells.Item(128, 2).Value = $sierraLeona[0]
$ws.Cells.Item(128, 3).Value = $sierraLeona[1]
$ws.Cells.Item(128, 4).Value = $sierraLeona[2]
$ws.Cells.Item(128, 5).Value = $sierraLeona[3]
$ws.Cells.Item(128, 6).Value = $sierraLeona[4]
$ws.Cells.Item(128, 7).Value = $sierraLeona[5]
$ws.Cells.Item(128, 8).Value = $sierraLeona[6]

# Row 129 -> Tanzania, keeping its previous figures
$ws.Cells.Item(129, 1).Value = "Tanzania"
$ws.Cells.Item(129, 2).Value = $tanzania[0]
$ws.Cells.Item(129, 3).Value = $tanzania[1]
$ws.Cells.Item(129, 4).Value = $tanzania[2]
$ws.Cells.Item(129, 5).Value = $tanzania[3]
$ws.Cells.Item(129, 6).Value = $tanzania[4]
$ws.Cells.Item(129, 7).Value = $tanzania[5]
$ws.Cells.Item(129, 8).Value = $tanzania[6]

# --- Refresh stats for Suazilandia (row 151) -------------------------------
$ws.Cells.Item(151, 2).Value = 205
$ws.Cells.Item(151, 3).Value = 2
$ws.Cells.Item(151, 4).Value = 78
$ws.Cells.Item(151, 5).Value = 125

# --- Refresh stats for Angola (row 178) ------------------------------------
$ws.Cells.Item(178, 2).Value = 50
$ws.Cells.Item(178, 3).Value = 2
$ws.Cells.Item(178, 5).Value = 30
$ws.Cells.Item(178, 7).Value = 1
$ws.Cells.Item(178, 8).Value = 3
